$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "H2-K1"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 322.4019826666667
$ws.Range("H2").Value = 967.205948
$ws.Range("I2").Value = 0.795674507658366
$ws.Range("J2").Value = 0.7956745076583662
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.179771666666667
$ws.Range("N2").Value = 6.539315
$ws.Range("O2").Value = 0.2349306639444428
$ws.Range("P2").Value = 0.2349306639444428
$ws.Range("Q2").Value = 702.7627070939579
$ws.Range("R2").Value = 6324.864363845621
$ws.Range("S2").Value = 0.1869283403678476
$ws.Range("T2").Value = 0.1869283403678476

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "H2-K1"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 322.4019826666667
$ws.Range("H3").Value = 967.205948
$ws.Range("I3").Value = 0.795674507658366
$ws.Range("J3").Value = 0.7956745076583662
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.237801
$ws.Range("N3").Value = 12.713403
$ws.Range("O3").Value = 0.4567402255103586
$ws.Range("P3").Value = 0.4567402255103586
$ws.Range("Q3").Value = 1366.275444546783
$ws.Range("R3").Value = 12296.47900092105
$ws.Range("S3").Value = 0.3634165540607256
$ws.Range("T3").Value = 0.3634165540607257

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "H2-K1"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 322.4019826666667
$ws.Range("H4").Value = 967.205948
$ws.Range("I4").Value = 0.795674507658366
$ws.Range("J4").Value = 0.7956745076583662
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.860789
$ws.Range("N4").Value = 8.582367000000001
$ws.Range("O4").Value = 0.3083291105451987
$ws.Range("P4").Value = 0.3083291105451986
$ws.Range("Q4").Value = 922.3240455909909
$ws.Range("R4").Value = 8300.916410318918
$ws.Range("S4").Value = 0.2453296132297929
$ws.Range("T4").Value = 0.2453296132297929

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "H2-K1"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 57.69151966666666
$ws.Range("H5").Value = 173.074559
$ws.Range("I5").Value = 0.1423802394983967
$ws.Range("J5").Value = 0.1423802394983967
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.179771666666667
$ws.Range("N5").Value = 6.539315
$ws.Range("O5").Value = 0.2349306639444428
$ws.Range("P5").Value = 0.2349306639444428
$ws.Range("Q5").Value = 125.7543399763428
$ws.Range("R5").Value = 1131.789059787085
$ws.Range("S5").Value = 0.03344948419792711
$ws.Range("T5").Value = 0.03344948419792711

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "H2-K1"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 57.69151966666666
$ws.Range("H6").Value = 173.074559
$ws.Range("I6").Value = 0.1423802394983967
$ws.Range("J6").Value = 0.1423802394983967
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.237801
$ws.Range("N6").Value = 12.713403
$ws.Range("O6").Value = 0.4567402255103586
$ws.Range("P6").Value = 0.4567402255103586
$ws.Range("Q6").Value = 244.4851797349197
$ws.Range("R6").Value = 2200.366617614277
$ws.Range("S6").Value = 0.06503078269671657
$ws.Range("T6").Value = 0.06503078269671658

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "H2-K1"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 57.69151966666666
$ws.Range("H7").Value = 173.074559
$ws.Range("I7").Value = 0.1423802394983967
$ws.Range("J7").Value = 0.1423802394983967
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.860789
$ws.Range("N7").Value = 8.582367000000001
$ws.Range("O7").Value = 0.3083291105451987
$ws.Range("P7").Value = 0.3083291105451986
$ws.Range("Q7").Value = 165.0432648556837
$ws.Range("R7").Value = 1485.389383701153
$ws.Range("S7").Value = 0.04389997260375302
$ws.Range("T7").Value = 0.04389997260375302

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "H2-K1"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 25.09980166666667
$ws.Range("H8").Value = 75.29940500000001
$ws.Range("I8").Value = 0.06194525284323717
$ws.Range("J8").Value = 0.06194525284323719
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.179771666666667
$ws.Range("N8").Value = 6.539315
$ws.Range("O8").Value = 0.2349306639444428
$ws.Range("P8").Value = 0.2349306639444428
$ws.Range("Q8").Value = 54.71183651195279
$ws.Range("R8").Value = 492.4065286075751
$ws.Range("S8").Value = 0.01455283937866809
$ws.Range("T8").Value = 0.01455283937866809

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "H2-K1"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 25.09980166666667
$ws.Range("H9").Value = 75.29940500000001
$ws.Range("I9").Value = 0.06194525284323717
$ws.Range("J9").Value = 0.06194525284323719
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237801
$ws.Range("N9").Value = 12.713403
$ws.Range("O9").Value = 0.4567402255103586
$ws.Range("P9").Value = 0.4567402255103586
$ws.Range("Q9").Value = 106.3679646028017
$ws.Range("R9").Value = 957.3116814252152
$ws.Range("S9").Value = 0.02829288875291633
$ws.Range("T9").Value = 0.02829288875291634

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "H2-K1"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 25.09980166666667
$ws.Range("H10").Value = 75.29940500000001
$ws.Range("I10").Value = 0.06194525284323717
$ws.Range("J10").Value = 0.06194525284323719
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.860789
$ws.Range("N10").Value = 8.582367000000001
$ws.Range("O10").Value = 0.3083291105451987
$ws.Range("P10").Value = 0.3083291105451986
$ws.Range("Q10").Value = 71.80523651018169
$ws.Range("R10").Value = 646.2471285916351
$ws.Range("S10").Value = 0.01909952471165276
$ws.Range("T10").Value = 0.01909952471165276
